$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7").Value = 15.058477850184744
$ws.Range("D7").Value = 138.6717582991117
$ws.Range("B8").Value = 931.5472412119683
$ws.Range("C8").Value = 6.7875746651481546
$ws.Range("D8").Value = 29.839568243479732
$ws.Range("B11").Value = 397.70676373500299
$ws.Range("D11").Value = 301.38197300776301
$ws.Range("C12").Value = 4.0878760483140084
$ws.Range("B14").Value = 314.00566328762528
$ws.Range("C14").Value = 4.6754703540978557
$ws.Range("D14").Value = 273.3793592273845
$ws.Range("B17").Value = 538.45403898299332
$ws.Range("C17").Value = 6.384756786916828
$ws.Range("D17").Value = 147.04411668022172
$ws.Range("B18").Value = 695.66648675853958
$ws.Range("C18").Value = 5.7005576276700198
$ws.Range("D18").Value = 257.73592117725241
$ws.Range("B19").Value = 327.4277695687523
$ws.Range("C19").Value = 5.4641980804167591
$ws.Range("D19").Value = 307.9228474015506
$ws.Range("B20").Value = 976.01382736066648
$ws.Range("C20").Value = 7.6689062639641694
$ws.Range("D20").Value = 35.092278793694753
$ws.Range("B21").Value = 436.60554443568412
$ws.Range("D21").Value = 235.09529315767605
$ws.Range("C22").Value = 5.0621892418420513
$ws.Range("D22").Value = 296.68077083850073
$ws.Range("B23").Value = 299.3075272149797
$ws.Range("C23").Value = 20.090947726133155
$ws.Range("D23").Value = 308.5817064150242
$ws.Range("B24").Value = 307.7470086895645
$ws.Range("D24").Value = 325.62832217178021
$ws.Range("B25").Value = 304.46849709606175
$ws.Range("D25").Value = 279.4059583612929
$ws.Range("B26").Value = 575.64301542973499
$ws.Range("C26").Value = 6.5173519677279899
$ws.Range("D26").Value = 34.884575505260514
$ws.Range("B27").Value = 1013.0934894050981
$ws.Range("C27").Value = 7.7199524064238139
$ws.Range("D27").Value = 33.770465267492625
$ws.Range("B28").Value = 1365.1523603800533
$ws.Range("C28").Value = 11.890056834598717
$ws.Range("D28").Value = 49.740982958699519
